$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1.Cells.Item(13, 8).Value = 9549.667
$ws1.Cells.Item(13, 9).Value = 3500
$ws1.Cells.Item(13, 10).Value = 10759.6
$ws1.Cells.Item(13, 11).Value = 3500
$ws1.Cells.Item(13, 12).Value = 10759.6
$ws1.Cells.Item(13, 13).Value = -3331
$ws1.Cells.Item(13, 14).Value = -11097.6

$ws1.Cells.Item(70, 8).Value = 900.2414
$ws1.Cells.Item(70, 10).Value = 793.8493
$ws1.Cells.Item(70, 12).Value = 2381.5479
$ws1.Cells.Item(70, 14).Value = -2921.5479

$ws1.Cells.Item(73, 8).Value = 900.2414
$ws1.Cells.Item(73, 10).Value = 793.8493
$ws1.Cells.Item(73, 12).Value = 2381.5479
$ws1.Cells.Item(73, 14).Value = -4253.5479

$ws1.Cells.Item(86, 8).Value = 6884.5386
$ws1.Cells.Item(86, 9).Value = 1999.6666
$ws1.Cells.Item(86, 11).Value = 1999.6666
$ws1.Cells.Item(86, 13).Value = -876.6666

$ws1.Cells.Item(89, 8).Value = 6884.5386
$ws1.Cells.Item(89, 9).Value = 1999.6666
$ws1.Cells.Item(89, 11).Value = 9998.333
$ws1.Cells.Item(89, 13).Value = -4382.333000000001

$ws1.Cells.Item(96, 8).Value = 1857.9
$ws1.Cells.Item(96, 9).Value = 1278
$ws1.Cells.Item(96, 10).Value = 3597.6
$ws1.Cells.Item(96, 11).Value = 3834
$ws1.Cells.Item(96, 12).Value = 10792.8
$ws1.Cells.Item(96, 13).Value = -2461
$ws1.Cells.Item(96, 14).Value = -13538.8

$ws1.Cells.Item(98, 8).Value = 595.7
$ws1.Cells.Item(98, 9).Value = 300.77777
$ws1.Cells.Item(98, 11).Value = 300.77777
$ws1.Cells.Item(98, 13).Value = 1197.22223

$ws1.Cells.Item(113, 8).Value = 7060.375
$ws1.Cells.Item(113, 9).Value = 4399.6665
$ws1.Cells.Item(113, 10).Value = 8656.8
$ws1.Cells.Item(113, 11).Value = 4399.6665
$ws1.Cells.Item(113, 12).Value = 8656.8
$ws1.Cells.Item(113, 13).Value = -1145.6665
$ws1.Cells.Item(113, 14).Value = -15164.8

$ws1.Cells.Item(122, 8).Value = 595.7
$ws1.Cells.Item(122, 9).Value = 300.77777
$ws1.Cells.Item(122, 11).Value = 902.33331
$ws1.Cells.Item(122, 13).Value = 1547.66669

$ws1.Cells.Item(132, 8).Value = 3743.7778
$ws1.Cells.Item(132, 9).Value = 3670.5715
$ws1.Cells.Item(132, 11).Value = 11011.7145
$ws1.Cells.Item(132, 13).Value = -8481.7145

$ws1.Cells.Item(135, 8).Value = 2062.8572
$ws1.Cells.Item(135, 9).Value = 576.8889
$ws1.Cells.Item(135, 10).Value = 4737.6
$ws1.Cells.Item(135, 11).Value = 5192.0001
$ws1.Cells.Item(135, 12).Value = 42638.4
$ws1.Cells.Item(135, 13).Value = -2657.0001
$ws1.Cells.Item(135, 14).Value = -47708.4

$ws1.Cells.Item(137, 8).Value = 2223.3513
$ws1.Cells.Item(137, 9).Value = 1717.5714
$ws1.Cells.Item(137, 10).Value = 3796.889
$ws1.Cells.Item(137, 11).Value = 5152.7142
$ws1.Cells.Item(137, 12).Value = 11390.667
$ws1.Cells.Item(137, 13).Value = -2602.7142
$ws1.Cells.Item(137, 14).Value = -16490.667

$ws1.Cells.Item(138, 8).Value = 4954
$ws1.Cells.Item(138, 9).Value = 3788.0833
$ws1.Cells.Item(138, 11).Value = 11364.2499
$ws1.Cells.Item(138, 13).Value = -6224.249899999999

$ws1.Cells.Item(140, 8).Value = 94394
$ws1.Cells.Item(140, 10).Value = 94394
$ws1.Cells.Item(140, 12).Value = 94394
$ws1.Cells.Item(140, 14).Value = -104754

$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Cells.Item(2, 8).Value = 4545.857
$ws2.Cells.Item(2, 9).Value = 4503.5
$ws2.Cells.Item(2, 11).Value = 4503.5
$ws2.Cells.Item(2, 13).Value = -4390.5

$ws2.Cells.Item(17, 8).Value = 44899
$ws2.Cells.Item(17, 10).Value = 44899
$ws2.Cells.Item(17, 12).Value = 44899
$ws2.Cells.Item(17, 14).Value = -45245

$ws2.Cells.Item(32, 8).Value = 6351.8804
$ws2.Cells.Item(32, 9).Value = 2616.3115
$ws2.Cells.Item(32, 11).Value = 2616.3115
$ws2.Cells.Item(32, 13).Value = -2329.3115

$ws2.Cells.Item(53, 8).Value = 34999
$ws2.Cells.Item(53, 9).Value = 29999
$ws2.Cells.Item(53, 11).Value = 29999
$ws2.Cells.Item(53, 13).Value = -29317

$ws2.Cells.Item(61, 8).Value = 35718090
$ws2.Cells.Item(61, 9).Value = 1564.091
$ws2.Cells.Item(61, 10).Value = 166678670
$ws2.Cells.Item(61, 11).Value = 1564.091
$ws2.Cells.Item(61, 12).Value = 166678670
$ws2.Cells.Item(61, 13).Value = -1352.091
$ws2.Cells.Item(61, 14).Value = -166679094

$ws2.Cells.Item(116, 8).Value = 4545.857
$ws2.Cells.Item(116, 9).Value = 4503.5
$ws2.Cells.Item(116, 11).Value = 4503.5
$ws2.Cells.Item(116, 13).Value = -2209.5

$ws2.Cells.Item(122, 8).Value = 7826.8086
$ws2.Cells.Item(122, 9).Value = 7843.628
$ws2.Cells.Item(122, 10).Value = 7646
$ws2.Cells.Item(122, 11).Value = 23530.884
$ws2.Cells.Item(122, 12).Value = 22938
$ws2.Cells.Item(122, 13).Value = -21080.884
$ws2.Cells.Item(122, 14).Value = -27838

$ws2.Cells.Item(136, 8).Value = 35718090
$ws2.Cells.Item(136, 9).Value = 1564.091
$ws2.Cells.Item(136, 10).Value = 166678670
$ws2.Cells.Item(136, 11).Value = 4692.272999999999
$ws2.Cells.Item(136, 12).Value = 500036010
$ws2.Cells.Item(136, 13).Value = -2142.272999999999
$ws2.Cells.Item(136, 14).Value = -500041110

$ws3 = $wb.Worksheets.Item("BSM")
$ws3.Cells.Item(3, 8).Value = 4545.857
$ws3.Cells.Item(3, 9).Value = 4503.5
$ws3.Cells.Item(3, 11).Value = 4503.5
$ws3.Cells.Item(3, 13).Value = -4389.5

$ws3.Cells.Item(81, 8).Value = 10999
$ws3.Cells.Item(81, 10).Value = 10999
$ws3.Cells.Item(81, 12).Value = 10999
$ws3.Cells.Item(81, 14).Value = -13121

$ws3.Cells.Item(84, 8).Value = 10999
$ws3.Cells.Item(84, 10).Value = 10999
$ws3.Cells.Item(84, 12).Value = 32997
$ws3.Cells.Item(84, 14).Value = -43605

$ws3.Cells.Item(86, 8).Value = 2073.6924
$ws3.Cells.Item(86, 9).Value = 1996.5
$ws3.Cells.Item(86, 11).Value = 1996.5
$ws3.Cells.Item(86, 13).Value = -873.5

$ws3.Cells.Item(89, 8).Value = 2073.6924
$ws3.Cells.Item(89, 9).Value = 1996.5
$ws3.Cells.Item(89, 11).Value = 9982.5
$ws3.Cells.Item(89, 13).Value = -4366.5

$ws3.Cells.Item(134, 8).Value = 6960.5
$ws3.Cells.Item(134, 9).Value = 1687.5
$ws3.Cells.Item(134, 11).Value = 5062.5
$ws3.Cells.Item(134, 13).Value = -2527.5

$ws4 = $wb.Worksheets.Item("CRP")
$ws4.Cells.Item(16, 8).Value = 814.0455
$ws4.Cells.Item(16, 9).Value = 537
$ws4.Cells.Item(16, 10).Value = 2568.6667
$ws4.Cells.Item(16, 11).Value = 537
$ws4.Cells.Item(16, 12).Value = 2568.6667
$ws4.Cells.Item(16, 13).Value = -250
$ws4.Cells.Item(16, 14).Value = -3142.6667

$ws4.Cells.Item(22, 8).Value = 626.2727
$ws4.Cells.Item(22, 9).Value = 282.33334
$ws4.Cells.Item(22, 10).Value = 1039
$ws4.Cells.Item(22, 11).Value = 282.33334
$ws4.Cells.Item(22, 12).Value = 1039
$ws4.Cells.Item(22, 13).Value = 67.66665999999998
$ws4.Cells.Item(22, 14).Value = -1739

$ws4.Cells.Item(31, 8).Value = 5452.525
$ws4.Cells.Item(31, 9).Value = 2175.2856
$ws4.Cells.Item(31, 11).Value = 2175.2856
$ws4.Cells.Item(31, 13).Value = -1880.2856

$ws4.Cells.Item(34, 8).Value = 5452.525
$ws4.Cells.Item(34, 9).Value = 2175.2856
$ws4.Cells.Item(34, 11).Value = 2175.2856
$ws4.Cells.Item(34, 13).Value = -1973.2856

$ws4.Cells.Item(58, 8).Value = 3589.36
$ws4.Cells.Item(58, 9).Value = 3318.6667
$ws4.Cells.Item(58, 10).Value = 5010.5
$ws4.Cells.Item(58, 11).Value = 3318.6667
$ws4.Cells.Item(58, 12).Value = 5010.5
$ws4.Cells.Item(58, 13).Value = -3115.6667
$ws4.Cells.Item(58, 14).Value = -5416.5

$ws4.Cells.Item(63, 8).Value = 95192.305
$ws4.Cells.Item(63, 9).Value = 0
$ws4.Cells.Item(63, 10).Value = 95192.305
$ws4.Cells.Item(63, 11).Value = 0
$ws4.Cells.Item(63, 12).Value = 95192.305
$ws4.Cells.Item(63, 13).Value = $null
$ws4.Cells.Item(63, 14).Value = -96564.305

$ws4.Cells.Item(66, 8).Value = 95192.305
$ws4.Cells.Item(66, 9).Value = 0
$ws4.Cells.Item(66, 10).Value = 95192.305
$ws4.Cells.Item(66, 11).Value = 0
$ws4.Cells.Item(66, 12).Value = 285576.915
$ws4.Cells.Item(66, 13).Value = $null
$ws4.Cells.Item(66, 14).Value = -292440.915

$ws4.Cells.Item(86, 8).Value = 5744.6665
$ws4.Cells.Item(86, 10).Value = 5744.6665
$ws4.Cells.Item(86, 12).Value = 5744.6665
$ws4.Cells.Item(86, 14).Value = -7990.6665

$ws4.Cells.Item(89, 8).Value = 5744.6665
$ws4.Cells.Item(89, 10).Value = 5744.6665
$ws4.Cells.Item(89, 12).Value = 28723.3325
$ws4.Cells.Item(89, 14).Value = -39955.3325

$ws4.Cells.Item(105, 8).Value = 200002980
$ws4.Cells.Item(105, 10).Value = 1000000000
$ws4.Cells.Item(105, 12).Value = 1000000000
$ws4.Cells.Item(105, 14).Value = -1000003494

$ws4.Cells.Item(107, 8).Value = 817.2308
$ws4.Cells.Item(107, 9).Value = 817.2308
$ws4.Cells.Item(107, 10).Value = 0
$ws4.Cells.Item(107, 11).Value = 817.2308
$ws4.Cells.Item(107, 12).Value = 0
$ws4.Cells.Item(107, 13).Value = 1102.7692
$ws4.Cells.Item(107, 14).Value = $null

$ws4.Cells.Item(113, 8).Value = 814.0455
$ws4.Cells.Item(113, 9).Value = 537
$ws4.Cells.Item(113, 10).Value = 2568.6667
$ws4.Cells.Item(113, 11).Value = 537
$ws4.Cells.Item(113, 12).Value = 2568.6667
$ws4.Cells.Item(113, 13).Value = 1633
$ws4.Cells.Item(113, 14).Value = -6908.6667

$ws4.Cells.Item(132, 8).Value = 5552.567
$ws4.Cells.Item(132, 9).Value = 3003.75
$ws4.Cells.Item(132, 11).Value = 9011.25
$ws4.Cells.Item(132, 13).Value = -6481.25

$ws4.Cells.Item(134, 8).Value = 5010.5454
$ws4.Cells.Item(134, 9).Value = 3911.8333
$ws4.Cells.Item(134, 10).Value = 15997.667
$ws4.Cells.Item(134, 11).Value = 11735.4999
$ws4.Cells.Item(134, 12).Value = 47993.001
$ws4.Cells.Item(134, 13).Value = -9200.499899999999
$ws4.Cells.Item(134, 14).Value = -53063.001

$ws4.Cells.Item(136, 8).Value = 3589.36
$ws4.Cells.Item(136, 9).Value = 3318.6667
$ws4.Cells.Item(136, 10).Value = 5010.5
$ws4.Cells.Item(136, 11).Value = 9956.000100000001
$ws4.Cells.Item(136, 12).Value = 15031.5
$ws4.Cells.Item(136, 13).Value = -7406.000100000001
$ws4.Cells.Item(136, 14).Value = -20131.5

$ws5 = $wb.Worksheets.Item("CUL")
$ws5.Cells.Item(3, 8).Value = 2400
$ws5.Cells.Item(3, 9).Value = 2400
$ws5.Cells.Item(3, 11).Value = 7200
$ws5.Cells.Item(3, 13).Value = -7088

$ws5.Cells.Item(5, 8).Value = 866.6667
$ws5.Cells.Item(5, 9).Value = 600
$ws5.Cells.Item(5, 10).Value = 1000
$ws5.Cells.Item(5, 11).Value = 1800
$ws5.Cells.Item(5, 12).Value = 3000
$ws5.Cells.Item(5, 13).Value = -1688
$ws5.Cells.Item(5, 14).Value = -3224

$ws5.Cells.Item(19, 8).Value = 1729.3334
$ws5.Cells.Item(19, 10).Value = 2295
$ws5.Cells.Item(19, 12).Value = 6885
$ws5.Cells.Item(19, 14).Value = -7233

$ws5.Cells.Item(48, 8).Value = 90
$ws5.Cells.Item(48, 9).Value = 90
$ws5.Cells.Item(48, 11).Value = 270
$ws5.Cells.Item(48, 13).Value = -20

$ws5.Cells.Item(49, 8).Value = 1533.3334
$ws5.Cells.Item(49, 9).Value = 1533.3334
$ws5.Cells.Item(49, 11).Value = 4600.0002
$ws5.Cells.Item(49, 13).Value = -4444.0002

$ws5.Cells.Item(51, 8).Value = 170
$ws5.Cells.Item(51, 10).Value = 150
$ws5.Cells.Item(51, 12).Value = 450
$ws5.Cells.Item(51, 14).Value = -1370

$ws5.Cells.Item(64, 8).Value = 4824.625
$ws5.Cells.Item(64, 9).Value = 4824.625
$ws5.Cells.Item(64, 11).Value = 14473.875
$ws5.Cells.Item(64, 13).Value = -14203.875

$ws5.Cells.Item(67, 8).Value = 4824.625
$ws5.Cells.Item(67, 9).Value = 4824.625
$ws5.Cells.Item(67, 11).Value = 14473.875
$ws5.Cells.Item(67, 13).Value = -13537.875

$ws5.Cells.Item(70, 8).Value = 14475.6
$ws5.Cells.Item(70, 9).Value = 10792.667
$ws5.Cells.Item(70, 10).Value = 20000
$ws5.Cells.Item(70, 11).Value = 32378.001
$ws5.Cells.Item(70, 12).Value = 60000
$ws5.Cells.Item(70, 13).Value = -32063.001
$ws5.Cells.Item(70, 14).Value = -60630

$ws5.Cells.Item(73, 8).Value = 14475.6
$ws5.Cells.Item(73, 9).Value = 10792.667
$ws5.Cells.Item(73, 10).Value = 20000
$ws5.Cells.Item(73, 11).Value = 32378.001
$ws5.Cells.Item(73, 12).Value = 60000
$ws5.Cells.Item(73, 13).Value = -31286.001
$ws5.Cells.Item(73, 14).Value = -62184

$ws5.Cells.Item(86, 8).Value = 227.14285
$ws5.Cells.Item(86, 9).Value = 198
$ws5.Cells.Item(86, 10).Value = 249
$ws5.Cells.Item(86, 11).Value = 594
$ws5.Cells.Item(86, 12).Value = 747
$ws5.Cells.Item(86, 13).Value = 592
$ws5.Cells.Item(86, 14).Value = -3119

$ws5.Cells.Item(89, 8).Value = 227.14285
$ws5.Cells.Item(89, 9).Value = 198
$ws5.Cells.Item(89, 10).Value = 249
$ws5.Cells.Item(89, 11).Value = 1782
$ws5.Cells.Item(89, 12).Value = 2241
$ws5.Cells.Item(89, 13).Value = 4146
$ws5.Cells.Item(89, 14).Value = -14097

$ws5.Cells.Item(103, 8).Value = 510.83334
$ws5.Cells.Item(103, 9).Value = 121.666664
$ws5.Cells.Item(103, 10).Value = 900
$ws5.Cells.Item(103, 11).Value = 364.999992
$ws5.Cells.Item(103, 12).Value = 2700
$ws5.Cells.Item(103, 13).Value = 514.000008
$ws5.Cells.Item(103, 14).Value = -4458

$ws5.Cells.Item(104, 8).Value = 423.83334
$ws5.Cells.Item(104, 9).Value = 373.5
$ws5.Cells.Item(104, 10).Value = 524.5
$ws5.Cells.Item(104, 11).Value = 1120.5
$ws5.Cells.Item(104, 12).Value = 1573.5
$ws5.Cells.Item(104, 13).Value = 1500.5
$ws5.Cells.Item(104, 14).Value = -6815.5

$ws5.Cells.Item(107, 8).Value = 269.42856
$ws5.Cells.Item(107, 10).Value = 269.42856
$ws5.Cells.Item(107, 12).Value = 808.28568
$ws5.Cells.Item(107, 14).Value = -4648.28568

$ws5.Cells.Item(109, 8).Value = 1000
$ws5.Cells.Item(109, 9).Value = 1000
$ws5.Cells.Item(109, 11).Value = 3000
$ws5.Cells.Item(109, 13).Value = -1960

$ws5.Cells.Item(112, 8).Value = 14321.368
$ws5.Cells.Item(112, 9).Value = 11939.556
$ws5.Cells.Item(112, 10).Value = 16465
$ws5.Cells.Item(112, 11).Value = 35818.66800000001
$ws5.Cells.Item(112, 12).Value = 49395
$ws5.Cells.Item(112, 13).Value = -34710.66800000001
$ws5.Cells.Item(112, 14).Value = -51611

$ws5.Cells.Item(121, 8).Value = 974230.75
$ws5.Cells.Item(121, 9).Value = 547.3
$ws5.Cells.Item(121, 10).Value = 4219842.5
$ws5.Cells.Item(121, 11).Value = 1641.9
$ws5.Cells.Item(121, 12).Value = 12659527.5
$ws5.Cells.Item(121, 13).Value = -331.8999999999999
$ws5.Cells.Item(121, 14).Value = -12662147.5

$ws5.Cells.Item(122, 8).Value = 583.46155
$ws5.Cells.Item(122, 9).Value = 587.44446
$ws5.Cells.Item(122, 10).Value = 574.5
$ws5.Cells.Item(122, 11).Value = 5287.00014
$ws5.Cells.Item(122, 12).Value = 5170.5
$ws5.Cells.Item(122, 13).Value = -2837.00014
$ws5.Cells.Item(122, 14).Value = -10070.5

$ws5.Cells.Item(128, 8).Value = 419994.2
$ws5.Cells.Item(128, 9).Value = 419994.2
$ws5.Cells.Item(128, 11).Value = 1259982.6
$ws5.Cells.Item(128, 13).Value = -1255002.6

$ws5.Cells.Item(133, 8).Value = 9300.65
$ws5.Cells.Item(133, 9).Value = 4008.0715
$ws5.Cells.Item(133, 11).Value = 12024.2145
$ws5.Cells.Item(133, 13).Value = -6964.2145

$ws5.Cells.Item(135, 8).Value = 866.6667
$ws5.Cells.Item(135, 9).Value = 600
$ws5.Cells.Item(135, 10).Value = 1000
$ws5.Cells.Item(135, 11).Value = 5400
$ws5.Cells.Item(135, 12).Value = 9000
$ws5.Cells.Item(135, 13).Value = -2865
$ws5.Cells.Item(135, 14).Value = -14070

$ws5.Cells.Item(137, 8).Value = 3926.75
$ws5.Cells.Item(137, 9).Value = 2468.6
$ws5.Cells.Item(137, 10).Value = 4968.2856
$ws5.Cells.Item(137, 11).Value = 7405.799999999999
$ws5.Cells.Item(137, 12).Value = 14904.8568
$ws5.Cells.Item(137, 13).Value = -2305.799999999999
$ws5.Cells.Item(137, 14).Value = -25104.8568

$ws6 = $wb.Worksheets.Item("GSM")
$ws6.Cells.Item(3, 8).Value = 11569.318
$ws6.Cells.Item(3, 9).Value = 15171.23
$ws6.Cells.Item(3, 10).Value = 6366.5557
$ws6.Cells.Item(3, 11).Value = 15171.23
$ws6.Cells.Item(3, 12).Value = 6366.5557
$ws6.Cells.Item(3, 13).Value = -15055.23
$ws6.Cells.Item(3, 14).Value = -6598.5557

$ws6.Cells.Item(6, 8).Value = 16846.5
$ws6.Cells.Item(6, 10).Value = 16846.5
$ws6.Cells.Item(6, 12).Value = 16846.5
$ws6.Cells.Item(6, 14).Value = -17072.5

$ws6.Cells.Item(9, 8).Value = 2389.2856
$ws6.Cells.Item(9, 9).Value = 444.33334
$ws6.Cells.Item(9, 10).Value = 3848
$ws6.Cells.Item(9, 11).Value = 444.33334
$ws6.Cells.Item(9, 12).Value = 3848
$ws6.Cells.Item(9, 13).Value = -274.33334
$ws6.Cells.Item(9, 14).Value = -4188

$ws6.Cells.Item(10, 8).Value = 47849.75
$ws6.Cells.Item(10, 9).Value = 10750
$ws6.Cells.Item(10, 11).Value = 10750
$ws6.Cells.Item(10, 13).Value = -10581

$ws6.Cells.Item(16, 8).Value = 16846.5
$ws6.Cells.Item(16, 10).Value = 16846.5
$ws6.Cells.Item(16, 12).Value = 16846.5
$ws6.Cells.Item(16, 14).Value = -17346.5

$ws6.Cells.Item(21, 8).Value = 10533.333
$ws6.Cells.Item(21, 9).Value = 10650
$ws6.Cells.Item(21, 10).Value = 9950
$ws6.Cells.Item(21, 11).Value = 10650
$ws6.Cells.Item(21, 12).Value = 9950
$ws6.Cells.Item(21, 13).Value = -10477
$ws6.Cells.Item(21, 14).Value = -10296

$ws6.Cells.Item(30, 8).Value = 10533.333
$ws6.Cells.Item(30, 9).Value = 10650
$ws6.Cells.Item(30, 10).Value = 9950
$ws6.Cells.Item(30, 11).Value = 10650
$ws6.Cells.Item(30, 12).Value = 9950
$ws6.Cells.Item(30, 13).Value = -10545
$ws6.Cells.Item(30, 14).Value = -10160

$ws6.Cells.Item(31, 8).Value = 6000
$ws6.Cells.Item(31, 9).Value = 6250
$ws6.Cells.Item(31, 11).Value = 6250
$ws6.Cells.Item(31, 13).Value = -5958

$ws6.Cells.Item(37, 8).Value = 6000
$ws6.Cells.Item(37, 9).Value = 6250
$ws6.Cells.Item(37, 11).Value = 6250
$ws6.Cells.Item(37, 13).Value = -5973

$ws6.Cells.Item(96, 8).Value = 199500
$ws6.Cells.Item(96, 10).Value = 199500
$ws6.Cells.Item(96, 12).Value = 199500
$ws6.Cells.Item(96, 14).Value = -204992

$ws6.Cells.Item(107, 8).Value = 496.1875
$ws6.Cells.Item(107, 9).Value = 328.25
$ws6.Cells.Item(107, 10).Value = 1000
$ws6.Cells.Item(107, 11).Value = 328.25
$ws6.Cells.Item(107, 12).Value = 1000
$ws6.Cells.Item(107, 13).Value = 1591.75
$ws6.Cells.Item(107, 14).Value = -4840

$ws6.Cells.Item(122, 8).Value = 6578.8486
$ws6.Cells.Item(122, 9).Value = 5454.5454
$ws6.Cells.Item(122, 10).Value = 8827.454
$ws6.Cells.Item(122, 11).Value = 16363.6362
$ws6.Cells.Item(122, 12).Value = 26482.362
$ws6.Cells.Item(122, 13).Value = -13913.6362
$ws6.Cells.Item(122, 14).Value = -31382.362

$ws7 = $wb.Worksheets.Item("LTW")
$ws7.Cells.Item(7, 8).Value = 5268.923
$ws7.Cells.Item(7, 10).Value = 5399.6
$ws7.Cells.Item(7, 14).Value = -5623.6

$ws7.Cells.Item(22, 8).Value = 3031.1428
$ws7.Cells.Item(22, 9).Value = 2857.3333
$ws7.Cells.Item(22, 10).Value = 3161.5
$ws7.Cells.Item(22, 11).Value = 2857.3333
$ws7.Cells.Item(22, 12).Value = 3161.5
$ws7.Cells.Item(22, 13).Value = -2562.3333
$ws7.Cells.Item(22, 14).Value = -3751.5

$ws7.Cells.Item(27, 8).Value = 3031.1428
$ws7.Cells.Item(27, 9).Value = 2857.3333
$ws7.Cells.Item(27, 10).Value = 3161.5
$ws7.Cells.Item(27, 11).Value = 2857.3333
$ws7.Cells.Item(27, 12).Value = 3161.5
$ws7.Cells.Item(27, 13).Value = -2750.3333
$ws7.Cells.Item(27, 14).Value = -3375.5

$ws7.Cells.Item(39, 8).Value = 55000
$ws7.Cells.Item(39, 10).Value = 55000
$ws7.Cells.Item(39, 12).Value = 55000
$ws7.Cells.Item(39, 14).Value = -55920

$ws7.Cells.Item(40, 8).Value = 3627.1765
$ws7.Cells.Item(40, 9).Value = 2724.1333
$ws7.Cells.Item(40, 10).Value = 10400
$ws7.Cells.Item(40, 11).Value = 2724.1333
$ws7.Cells.Item(40, 12).Value = 10400
$ws7.Cells.Item(40, 13).Value = -2588.1333
$ws7.Cells.Item(40, 14).Value = -10672

$ws7.Cells.Item(45, 8).Value = 21874.75
$ws7.Cells.Item(45, 9).Value = 8750
$ws7.Cells.Item(45, 11).Value = 8750
$ws7.Cells.Item(45, 13).Value = -8343

$ws7.Cells.Item(46, 8).Value = 4500.5
$ws7.Cells.Item(46, 9).Value = 1999.5
$ws7.Cells.Item(46, 11).Value = 1999.5
$ws7.Cells.Item(46, 13).Value = -1811.5

$ws7.Cells.Item(48, 8).Value = 49999
$ws7.Cells.Item(48, 9).Value = 0
$ws7.Cells.Item(48, 10).Value = 49999
$ws7.Cells.Item(48, 11).Value = 0
$ws7.Cells.Item(48, 12).Value = 49999
$ws7.Cells.Item(48, 13).Value = $null
$ws7.Cells.Item(48, 14).Value = -51321

$ws7.Cells.Item(55, 8).Value = 1320.6666
$ws7.Cells.Item(55, 9).Value = 271.2
$ws7.Cells.Item(55, 10).Value = 2070.2856
$ws7.Cells.Item(55, 11).Value = 271.2
$ws7.Cells.Item(55, 12).Value = 2070.2856
$ws7.Cells.Item(55, 13).Value = -98.19999999999999
$ws7.Cells.Item(55, 14).Value = -2416.2856

$ws7.Cells.Item(61, 8).Value = 3178
$ws7.Cells.Item(61, 9).Value = 916.6667
$ws7.Cells.Item(61, 10).Value = 6570
$ws7.Cells.Item(61, 11).Value = 916.6667
$ws7.Cells.Item(61, 12).Value = 6570
$ws7.Cells.Item(61, 13).Value = -714.6667
$ws7.Cells.Item(61, 14).Value = -6974

$ws7.Cells.Item(68, 8).Value = 7035.7856
$ws7.Cells.Item(68, 9).Value = 7677.778
$ws7.Cells.Item(68, 10).Value = 5880.2
$ws7.Cells.Item(68, 11).Value = 7677.778
$ws7.Cells.Item(68, 12).Value = 5880.2
$ws7.Cells.Item(68, 13).Value = -6928.778
$ws7.Cells.Item(68, 14).Value = -7378.2

$ws7.Cells.Item(71, 8).Value = 7035.7856
$ws7.Cells.Item(71, 9).Value = 7677.778
$ws7.Cells.Item(71, 10).Value = 5880.2
$ws7.Cells.Item(71, 11).Value = 38388.89
$ws7.Cells.Item(71, 12).Value = 29401
$ws7.Cells.Item(71, 13).Value = -34644.89
$ws7.Cells.Item(71, 14).Value = -36889

$ws7.Cells.Item(113, 8).Value = 3178
$ws7.Cells.Item(113, 9).Value = 916.6667
$ws7.Cells.Item(113, 10).Value = 6570
$ws7.Cells.Item(113, 11).Value = 916.6667
$ws7.Cells.Item(113, 12).Value = 6570
$ws7.Cells.Item(113, 13).Value = 1253.3333
$ws7.Cells.Item(113, 14).Value = -10910

$ws7.Cells.Item(122, 8).Value = 6203.3335
$ws7.Cells.Item(122, 9).Value = 5676.591
$ws7.Cells.Item(122, 10).Value = 11997.5
$ws7.Cells.Item(122, 11).Value = 17029.773
$ws7.Cells.Item(122, 12).Value = 35992.5
$ws7.Cells.Item(122, 13).Value = -14579.773
$ws7.Cells.Item(122, 14).Value = -40892.5

$ws7.Cells.Item(126, 8).Value = 5268.923
$ws7.Cells.Item(126, 10).Value = 5399.6
$ws7.Cells.Item(126, 12).Value = 16198.8
$ws7.Cells.Item(126, 14).Value = -21138.8

$ws7.Cells.Item(132, 8).Value = 11715.23
$ws7.Cells.Item(132, 9).Value = 9921.777
$ws7.Cells.Item(132, 11).Value = 29765.331
$ws7.Cells.Item(132, 13).Value = -27235.331

$ws7.Cells.Item(136, 8).Value = 16134686
$ws7.Cells.Item(136, 9).Value = 3458.862
$ws7.Cells.Item(136, 10).Value = 30310612
$ws7.Cells.Item(136, 11).Value = 10376.586
$ws7.Cells.Item(136, 12).Value = 90931836
$ws7.Cells.Item(136, 13).Value = -7826.585999999999
$ws7.Cells.Item(136, 14).Value = -90936936

$ws8 = $wb.Worksheets.Item("WVR")
$ws8.Cells.Item(9, 8).Value = 57450
$ws8.Cells.Item(9, 10).Value = 57450
$ws8.Cells.Item(9, 12).Value = 57450
$ws8.Cells.Item(9, 14).Value = -57730

$ws8.Cells.Item(12, 8).Value = 20000
$ws8.Cells.Item(12, 10).Value = 0
$ws8.Cells.Item(12, 12).Value = 0
$ws8.Cells.Item(12, 14).Value = $null

$ws8.Cells.Item(14, 8).Value = 20860.8
$ws8.Cells.Item(14, 9).Value = 34900
$ws8.Cells.Item(14, 10).Value = 6821.6
$ws8.Cells.Item(14, 11).Value = 34900
$ws8.Cells.Item(14, 12).Value = 6821.6
$ws8.Cells.Item(14, 13).Value = -34732
$ws8.Cells.Item(14, 14).Value = -7157.6

$ws8.Cells.Item(38, 8).Value = 50000
$ws8.Cells.Item(38, 10).Value = 0
$ws8.Cells.Item(38, 12).Value = 0
$ws8.Cells.Item(38, 14).Value = $null

$ws8.Cells.Item(42, 8).Value = 76633.336
$ws8.Cells.Item(42, 10).Value = 76633.336
$ws8.Cells.Item(42, 12).Value = 76633.336
$ws8.Cells.Item(42, 14).Value = -77389.336

$ws8.Cells.Item(51, 8).Value = 47500
$ws8.Cells.Item(51, 10).Value = 47500
$ws8.Cells.Item(51, 12).Value = 47500
$ws8.Cells.Item(51, 14).Value = -48520

$ws8.Cells.Item(52, 8).Value = 20000
$ws8.Cells.Item(52, 9).Value = 20000
$ws8.Cells.Item(52, 11).Value = 20000
$ws8.Cells.Item(52, 13).Value = -19774

$ws8.Cells.Item(62, 8).Value = 5818.8887
$ws8.Cells.Item(62, 9).Value = 5818.8887
$ws8.Cells.Item(62, 10).Value = 0
$ws8.Cells.Item(62, 11).Value = 5818.8887
$ws8.Cells.Item(62, 12).Value = 0
$ws8.Cells.Item(62, 13).Value = -5194.8887
$ws8.Cells.Item(62, 14).Value = $null

$ws8.Cells.Item(65, 8).Value = 5818.8887
$ws8.Cells.Item(65, 9).Value = 5818.8887
$ws8.Cells.Item(65, 10).Value = 0
$ws8.Cells.Item(65, 11).Value = 29094.4435
$ws8.Cells.Item(65, 12).Value = 0
$ws8.Cells.Item(65, 13).Value = -25974.4435
$ws8.Cells.Item(65, 14).Value = $null

$ws8.Cells.Item(81, 8).Value = 2640
$ws8.Cells.Item(81, 9).Value = 2460
$ws8.Cells.Item(81, 10).Value = 3000
$ws8.Cells.Item(81, 11).Value = 4920
$ws8.Cells.Item(81, 12).Value = 6000
$ws8.Cells.Item(81, 13).Value = -3859
$ws8.Cells.Item(81, 14).Value = -8122

$ws8.Cells.Item(84, 8).Value = 2640
$ws8.Cells.Item(84, 9).Value = 2460
$ws8.Cells.Item(84, 10).Value = 3000
$ws8.Cells.Item(84, 11).Value = 24600
$ws8.Cells.Item(84, 12).Value = 30000
$ws8.Cells.Item(84, 13).Value = -19296
$ws8.Cells.Item(84, 14).Value = -40608

$ws8.Cells.Item(96, 8).Value = 26052.111
$ws8.Cells.Item(96, 9).Value = 3000
$ws8.Cells.Item(96, 10).Value = 28933.625
$ws8.Cells.Item(96, 11).Value = 3000
$ws8.Cells.Item(96, 12).Value = 28933.625
$ws8.Cells.Item(96, 13).Value = -1627
$ws8.Cells.Item(96, 14).Value = -31679.625

$ws8.Cells.Item(100, 8).Value = 429.0909
$ws8.Cells.Item(100, 9).Value = 395.7143
$ws8.Cells.Item(100, 11).Value = 791.4286
$ws8.Cells.Item(100, 13).Value = -250.4286

$ws8.Cells.Item(113, 8).Value = 876.7778
$ws8.Cells.Item(113, 9).Value = 479.6
$ws8.Cells.Item(113, 11).Value = 1438.8
$ws8.Cells.Item(113, 13).Value = 731.1999999999998

$ws8.Cells.Item(122, 8).Value = 2818.125
$ws8.Cells.Item(122, 9).Value = 2706
$ws8.Cells.Item(122, 11).Value = 8118
$ws8.Cells.Item(122, 13).Value = -5668

$ws8.Cells.Item(126, 8).Value = 3112.375
$ws8.Cells.Item(126, 9).Value = 2899.8333
$ws8.Cells.Item(126, 11).Value = 8699.499899999999
$ws8.Cells.Item(126, 13).Value = -6229.499899999999

$ws8.Cells.Item(132, 8).Value = 2172.7693
$ws8.Cells.Item(132, 9).Value = 1520.7084
$ws8.Cells.Item(132, 10).Value = 9997.5
$ws8.Cells.Item(132, 11).Value = 4562.1252
$ws8.Cells.Item(132, 12).Value = 29992.5
$ws8.Cells.Item(132, 13).Value = -2032.1252
$ws8.Cells.Item(132, 14).Value = -35052.5

$ws8.Cells.Item(136, 8).Value = 4932.778
$ws8.Cells.Item(136, 9).Value = 4820.857
$ws8.Cells.Item(136, 10).Value = 5324.5
$ws8.Cells.Item(136, 11).Value = 14462.571
$ws8.Cells.Item(136, 12).Value = 15973.5
$ws8.Cells.Item(136, 13).Value = -11912.571
$ws8.Cells.Item(136, 14).Value = -21073.5

$ws8.Cells.Item(137, 8).Value = 82999
$ws8.Cells.Item(137, 9).Value = 0
$ws8.Cells.Item(137, 11).Value = 0
$ws8.Cells.Item(137, 13).Value = $null
